# Updates the cryptocurrency price/volume table on the active worksheet.
#
# Columns D (Price) and E (Volume(1h)) hold plain-text values in the source
# data (e.g. "63.034.71", "1.00", "  -4.89%  ") even though many of them look
# numeric. Two pairs of rows (42/43 and 47/48) also swap rank position, so
# their B (Coin) / C (Link) columns change too.
#
# Assigning a numeric-looking string straight to Range.Value lets Excel's
# COM layer "smart" re-type it as a real number (e.g. "1.00" -> 1, "6.36" ->
# 6.3600000000000003), which would corrupt the text. To avoid that we
# temporarily force Text format ("@") on each cell before writing, then
# ClearFormats() immediately after so the cell reverts to the workbook's
# original (unstyled) look while keeping the text content we just wrote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ D="63.199.16"; E="  -4.72%  " }
    3 = @{ D="3.084.42"; E="  -4.74%  " }
    4 = @{ D="1.00"; E="  +0.07%  " }
    5 = @{ D="549.98"; E="  -4.43%  " }
    6 = @{ D="136.03"; E="  -12.10%  " }
    7 = @{ E="  +0.11%  " }
    8 = @{ D="3.077.96"; E="  -4.66%  " }
    9 = @{ D="0.495"; E="  -3.94%  " }
    10 = @{ D="0.157"; E="  -6.34%  " }
    11 = @{ D="6.36"; E="  -10.42%  " }
    12 = @{ D="0.465"; E="  -4.25%  " }
    13 = @{ D="35.17"; E="  -7.60%  " }
    14 = @{ D="0.0000218"; E="  -7.53%  " }
    15 = @{ D="3.580.92"; E="  -4.53%  " }
    16 = @{ D="63.187.27"; E="  -4.60%  " }
    17 = @{ E="  -3.45%  " }
    18 = @{ D="3.073.31"; E="  -5.08%  " }
    19 = @{ D="6.70"; E="  -5.33%  " }
    20 = @{ D="491.16"; E="  -10.62%  " }
    21 = @{ D="13.56"; E="  -6.33%  " }
    22 = @{ D="0.711"; E="  -4.04%  " }
    23 = @{ D="7.22"; E="  -7.64%  " }
    24 = @{ D="78.29"; E="  -4.66%  " }
    25 = @{ D="12.37"; E="  -8.28%  " }
    26 = @{ D="1.00"; E="  +0.21%  " }
    27 = @{ D="8.46"; E="  -10.21%  " }
    28 = @{ D="2.75"; E="  -5.44%  " }
    29 = @{ D="0.998"; E="  -0.43%  " }
    30 = @{ D="1.97"; E="  -13.45%  " }
    31 = @{ D="26.50"; E="  -5.31%  " }
    32 = @{ E="  -3.28%  " }
    33 = @{ D="2.53"; E="  -8.70%  " }
    34 = @{ D="58.33"; E="  +6.53%  " }
    35 = @{ D="522.51"; E="  -8.41%  " }
    36 = @{ D="6.02"; E="  -7.02%  " }
    37 = @{ D="5.14"; E="  -11.04%  " }
    38 = @{ D="0.0409"; E="  -12.85%  " }
    39 = @{ D="3.105.20"; E="  -0.94%  " }
    40 = @{ D="0.0802"; E="  -8.55%  " }
    41 = @{ D="0.120"; E="  -5.32%  " }
    42 = @{ B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="2.68"; E="  -11.51%  " }
    43 = @{ B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="8.14"; E="  -5.94%  " }
    44 = @{ D="0.257"; E="  -6.29%  " }
    45 = @{ E="  +0.11%  " }
    46 = @{ D="2.10"; E="  -10.53%  " }
    47 = @{ B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="24.85"; E="  -8.70%  " }
    48 = @{ B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="121.24"; E="  -1.07%  " }
    49 = @{ E="  -4.82%  " }
    50 = @{ D=("0.0" + ([string][char]0x2083) + "0503"); E="  -11.56%  " }
    51 = @{ D="2.35"; E="  +49.80%  " }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col]
        $cell.ClearFormats()
    }
}

